# Fixed path in excel files: replace backslash-escaped sql paths with
# forward-slash paths, and update the active selection on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the source / target SQL query path cells to use forward slashes
# instead of escaped backslashes.
$ws.Range("B15").Value = "test/sql/sourcesql"
$ws.Range("B28").Value = "test/sql/targetsql"

# Update the saved selection/active cell on the sheet: it now shows the
# whole column B selected with B1 as the active cell (instead of a single
# selected cell B33).
$ws.Columns("B").Select()
